# Update "想去人数" (want-to-go count) values in the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 294
$ws1.Range("F3").Value = 1185
$ws1.Range("F4").Value = 16757
$ws1.Range("F6").Value = 1639
$ws1.Range("F7").Value = 62
$ws1.Range("F9").Value = 375
$ws1.Range("F10").Value = 214
$ws1.Range("F11").Value = 126
$ws1.Range("F12").Value = 11622
$ws1.Range("F14").Value = 1300
$ws1.Range("F15").Value = 4604
$ws1.Range("F16").Value = 429
$ws1.Range("F18").Value = 66
$ws1.Range("F19").Value = 888

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 294
$ws4.Range("F4").Value = 1185
$ws4.Range("F5").Value = 16757
$ws4.Range("F7").Value = 1639
$ws4.Range("F8").Value = 62
$ws4.Range("F10").Value = 375
$ws4.Range("F11").Value = 214
$ws4.Range("F12").Value = 126
$ws4.Range("F15").Value = 11622
$ws4.Range("F17").Value = 1301
$ws4.Range("F18").Value = 4604
$ws4.Range("F19").Value = 429
$ws4.Range("F21").Value = 66
$ws4.Range("F22").Value = 888
